$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.696.53"
$ws.Range("D3").Value = "1.947.66"
$ws.Range("E3").Value = "  +1.78%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'246.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.92%  "
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").Value = "'0.2939"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.63%  "
$ws.Range("D9").Value = "'0.06821"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.60%  "
$ws.Range("D10").Value = "'112.60"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.22%  "
$ws.Range("D11").Value = "'19.42"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.78%  "
$ws.Range("D12").Value = "1.951.36"
$ws.Range("E12").Value = "  +1.99%  "
$ws.Range("D13").Value = "'0.07673"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'5.518"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.57%  "
$ws.Range("D15").Value = "'0.6899"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.13%  "
$ws.Range("D16").Value = "'296.47"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.21%  "
$ws.Range("D17").Value = "30.735.35"
$ws.Range("E17").Value = "  +0.78%  "
$ws.Range("E18").Value = "  +3.30%  "
$ws.Range("D19").Value = "'5.676"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.57%  "
$ws.Range("D20").Value = "'0.000007697"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.14%  "
$ws.Range("D21").Value = "2.202.72"
$ws.Range("E21").Value = "  +1.83%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").Value = "'1.000"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").Value = "'6.603"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.60%  "
$ws.Range("D25").Value = "'9.872"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = "'168.26"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.81%  "
$ws.Range("E27").Value = "  +0.83%  "
$ws.Range("E28").Value = "  +3.46%  "
$ws.Range("D29").Value = "'0.1088"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.22%  "
$ws.Range("D30").Value = "'1.437"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.68%  "
$ws.Range("D31").Value = "'4.723"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +16.55%  "
$ws.Range("D32").Value = "'4.522"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.93%  "
$ws.Range("D33").Value = "'0.05094"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.09%  "
$ws.Range("D34").Value = "'0.7777"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.65%  "
$ws.Range("E35").Value = "  +2.71%  "
$ws.Range("D36").Value = "'0.02087"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.77%  "
$ws.Range("D37").Value = "'2.733"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").Value = "'2.701"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.25%  "
$ws.Range("D39").Value = "'2.057"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.95%  "
$ws.Range("D40").Value = "'111.52"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.64%  "
$ws.Range("D41").Value = "'0.4467"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.87%  "
$ws.Range("D42").Value = "'0.8737"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.02%  "
$ws.Range("D43").Value = "'5.921"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.64%  "
$ws.Range("D44").Value = "'70.22"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.43%  "
$ws.Range("E45").Value = "  +0.30%  "
$ws.Range("D46").Value = "'7.380"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.29%  "
$ws.Range("D47").Value = "'9.406"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.52%  "
$ws.Range("D48").Value = "'48.63"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.67%  "
$ws.Range("D49").Value = "'0.1257"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.65%  "
$ws.Range("D50").Value = "'35.80"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.08%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "'1.485"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.23%  "
